# Apply scheduled-runner profit recalculation updates to the Leve profit
# tracking sheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR). Values below are the
# refreshed currentAveragePrice / LevePrice / LeveProfit figures pulled in
# by the runner; a handful of rows also gain/lose their NQ profit cell
# (M) or HQ profit cell (N) when that side of the recipe is no longer priced.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 38
$ws.Range("H38").Value = 1477.1052
$ws.Range("I38").Value = 369.54544
$ws.Range("J38").Value = 3000
$ws.Range("K38").Value = 1108.63632
$ws.Range("L38").Value = 9000
$ws.Range("M38").Value = -736.6363200000001
$ws.Range("N38").Value = -9744

# row 40
$ws.Range("H40").Value = 62552.23
$ws.Range("I40").Value = 100810
$ws.Range("J40").Value = 1339.8
$ws.Range("K40").Value = 100810
$ws.Range("L40").Value = 1339.8
$ws.Range("M40").Value = -100635
$ws.Range("N40").Value = -1689.8

# row 58
$ws.Range("H58").Value = 19413.225
$ws.Range("J58").Value = 22459.14
$ws.Range("L58").Value = 67377.42
$ws.Range("N58").Value = -67677.42

# row 87
$ws.Range("H87").Value = 20965.95
$ws.Range("J87").Value = 20965.95
$ws.Range("L87").Value = 20965.95
$ws.Range("N87").Value = -23461.95

# row 90
$ws.Range("H90").Value = 20965.95
$ws.Range("J90").Value = 20965.95
$ws.Range("L90").Value = 62897.85000000001
$ws.Range("N90").Value = -75377.85000000001

# row 132
$ws.Range("H132").Value = 2565829.2
$ws.Range("I132").Value = 2779316.8
$ws.Range("J132").Value = 3980.3333
$ws.Range("K132").Value = 8337950.399999999
$ws.Range("L132").Value = 11940.9999
$ws.Range("M132").Value = -8335420.399999999
$ws.Range("N132").Value = -17000.9999

# row 138
$ws.Range("H138").Value = 2684.6736
$ws.Range("I138").Value = 941.4146
$ws.Range("J138").Value = 3938.5964
$ws.Range("K138").Value = 2824.2438
$ws.Range("L138").Value = 11815.7892
$ws.Range("M138").Value = 2315.7562
$ws.Range("N138").Value = -22095.7892

$ws = $wb.Worksheets.Item("ARM")
# row 44
$ws.Range("H44").Value = 19366
$ws.Range("J44").Value = 19366
$ws.Range("L44").Value = 19366
$ws.Range("N44").Value = -20342

# row 55
$ws.Range("H55").Value = 19926.5
$ws.Range("J55").Value = 19926.5
$ws.Range("L55").Value = 19926.5
$ws.Range("N55").Value = -20556.5

# row 74
$ws.Range("H74").Value = 509.72916
$ws.Range("I74").Value = 472.3111
$ws.Range("J74").Value = 1071
$ws.Range("K74").Value = 472.3111
$ws.Range("L74").Value = 1071
$ws.Range("M74").Value = 401.6889
$ws.Range("N74").Value = -2819

# row 76
$ws.Range("H76").Value = 27090.818
$ws.Range("J76").Value = 27090.818
$ws.Range("L76").Value = 27090.818
$ws.Range("N76").Value = -27766.818

# row 77
$ws.Range("H77").Value = 509.72916
$ws.Range("I77").Value = 472.3111
$ws.Range("J77").Value = 1071
$ws.Range("K77").Value = 2361.5555
$ws.Range("L77").Value = 5355
$ws.Range("M77").Value = 2006.4445
$ws.Range("N77").Value = -14091

# row 79
$ws.Range("H79").Value = 27090.818
$ws.Range("J79").Value = 27090.818
$ws.Range("L79").Value = 27090.818
$ws.Range("N79").Value = -29430.818

# row 80
$ws.Range("H80").Value = 27142.666
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 27142.666
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 27142.666
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -29138.666

# row 83
$ws.Range("H83").Value = 27142.666
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 27142.666
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 81427.99800000001
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -91411.99800000001

# row 110
$ws.Range("H110").Value = 1355.2069
$ws.Range("I110").Value = 670.9545000000001
$ws.Range("K110").Value = 670.9545000000001
$ws.Range("M110").Value = 1374.0455

# row 139
$ws.Range("H139").Value = 21386.924
$ws.Range("J139").Value = 21386.924
$ws.Range("L139").Value = 21386.924
$ws.Range("N139").Value = -31666.924

$ws = $wb.Worksheets.Item("CRP")
# row 22
$ws.Range("H22").Value = 5750
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 5750
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 5750
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -6450

# row 52
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

# row 70
$ws.Range("H70").Value = 34999.8
$ws.Range("J70").Value = 35000
$ws.Range("L70").Value = 35000
$ws.Range("N70").Value = -35630

# row 73
$ws.Range("H73").Value = 34999.8
$ws.Range("J73").Value = 35000
$ws.Range("L73").Value = 35000
$ws.Range("N73").Value = -37184

$ws = $wb.Worksheets.Item("CUL")
# row 17
$ws.Range("H17").Value = 485.53333
$ws.Range("I17").Value = 246.71428
$ws.Range("J17").Value = 694.5
$ws.Range("K17").Value = 740.14284
$ws.Range("L17").Value = 2083.5
$ws.Range("M17").Value = -571.14284
$ws.Range("N17").Value = -2421.5

# row 34
$ws.Range("H34").Value = 7393.5293
$ws.Range("I34").Value = 248.33333
$ws.Range("J34").Value = 11290.909
$ws.Range("K34").Value = 744.99999
$ws.Range("L34").Value = 33872.727
$ws.Range("M34").Value = -660.99999
$ws.Range("N34").Value = -34040.727

# row 39
$ws.Range("H39").Value = 1742.2
$ws.Range("I39").Value = 450
$ws.Range("J39").Value = 1885.7778
$ws.Range("K39").Value = 1350
$ws.Range("L39").Value = 5657.3334
$ws.Range("M39").Value = -1056
$ws.Range("N39").Value = -6245.3334

# row 55
$ws.Range("H55").Value = 2691.7407
$ws.Range("J55").Value = 3077.261
$ws.Range("L55").Value = 9231.782999999999
$ws.Range("N55").Value = -9585.782999999999

# row 134
$ws.Range("H134").Value = 2723
$ws.Range("I134").Value = 1076
$ws.Range("J134").Value = 4187
$ws.Range("K134").Value = 3228
$ws.Range("L134").Value = 12561
$ws.Range("M134").Value = 1842
$ws.Range("N134").Value = -22701

# row 140
$ws.Range("H140").Value = 10420876
$ws.Range("I140").Value = 33334074
$ws.Range("J140").Value = 5785.4546
$ws.Range("K140").Value = 100002222
$ws.Range("L140").Value = 17356.3638
$ws.Range("M140").Value = -99997042
$ws.Range("N140").Value = -27716.3638

$ws = $wb.Worksheets.Item("GSM")
# row 43
$ws.Range("H43").Value = 4378.5
$ws.Range("I43").Value = 2318
$ws.Range("K43").Value = 2318
$ws.Range("M43").Value = -2167

# row 46
$ws.Range("H46").Value = 8878.25
$ws.Range("I46").Value = 2000
$ws.Range("K46").Value = 2000
$ws.Range("M46").Value = -1844

# row 57
$ws.Range("H57").Value = 18061
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 18061
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 18061
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -19701

# row 80
$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 2928.5715
$ws.Range("J80").Value = 3250
$ws.Range("K80").Value = 2928.5715
$ws.Range("L80").Value = 3250
$ws.Range("M80").Value = -1930.5715
$ws.Range("N80").Value = -5246

# row 83
$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 2928.5715
$ws.Range("J83").Value = 3250
$ws.Range("K83").Value = 14642.8575
$ws.Range("L83").Value = 16250
$ws.Range("M83").Value = -9650.8575
$ws.Range("N83").Value = -26234

# row 109
$ws.Range("H109").Value = 10285
$ws.Range("J109").Value = 10285
$ws.Range("L109").Value = 10285
$ws.Range("N109").Value = -12365

# row 126
$ws.Range("H126").Value = 2901.394
$ws.Range("I126").Value = 1954.7693
$ws.Range("J126").Value = 3516.7
$ws.Range("K126").Value = 5864.3079
$ws.Range("L126").Value = 10550.1
$ws.Range("M126").Value = -3394.3079
$ws.Range("N126").Value = -15490.1

$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 200002530
$ws.Range("I22").Value = 500000350
$ws.Range("J22").Value = 3966.6667
$ws.Range("K22").Value = 500000350
$ws.Range("L22").Value = 3966.6667
$ws.Range("M22").Value = -500000055
$ws.Range("N22").Value = -4556.6667

# row 27
$ws.Range("H27").Value = 200002530
$ws.Range("I27").Value = 500000350
$ws.Range("J27").Value = 3966.6667
$ws.Range("K27").Value = 500000350
$ws.Range("L27").Value = 3966.6667
$ws.Range("M27").Value = -500000243
$ws.Range("N27").Value = -4180.6667

# row 46
$ws.Range("H46").Value = 1627.7778
$ws.Range("J46").Value = 2069.2307
$ws.Range("L46").Value = 2069.2307
$ws.Range("N46").Value = -2445.2307

# row 55
$ws.Range("H55").Value = 1317.5454
$ws.Range("I55").Value = 225.25
$ws.Range("K55").Value = 225.25
$ws.Range("M55").Value = -52.25

$ws = $wb.Worksheets.Item("WVR")
# row 3
$ws.Range("H3").Value = 1720225.8
$ws.Range("I3").Value = 2400116
$ws.Range("J3").Value = 20500
$ws.Range("K3").Value = 2400116
$ws.Range("L3").Value = 20500
$ws.Range("M3").Value = -2400002
$ws.Range("N3").Value = -20728

# row 54
$ws.Range("H54").Value = 13729.75
$ws.Range("J54").Value = 13729.75
$ws.Range("L54").Value = 13729.75
$ws.Range("N54").Value = -14769.75

# row 81
$ws.Range("H81").Value = 888.2857
$ws.Range("I81").Value = 888.2857
$ws.Range("K81").Value = 1776.5714
$ws.Range("M81").Value = -715.5714

# row 84
$ws.Range("H84").Value = 888.2857
$ws.Range("I84").Value = 888.2857
$ws.Range("K84").Value = 8882.857
$ws.Range("M84").Value = -3578.857
